$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at row 6 (new July 2025 month at the top of the table,
# below the header row), pushing the existing data down by one row.
# -4162 = xlFormatFromRightOrBelow so the new row's formatting matches the
# row that ends up just below it (old row 6 -> new row 7's sibling style).
$ws.Rows.Item(6).Insert(-4162)

# Copy formatting from what is now row 8 (the "odd" banded style) onto the
# freshly inserted row 6 so the visual banding pattern is preserved exactly
# as authored.
$ws.Range("B8:G8").Copy()
$ws.Range("B6:G6").PasteSpecial(-4122)

# Populate the new row with the July 2025 figures.
$ws.Range("B6").Value = 2025
$ws.Range("C6").Value = "Jul."
$ws.Range("D6").Value = 365.85
$ws.Range("E6").Value = 34632.449
$ws.Range("F6").Value = 4369.473
$ws.Range("G6").Value = 111.468

# Grow the "Tabla1" table/autofilter range by one row so it keeps covering
# the whole dataset (B5:G95 -> B5:G96).
$t = $ws.ListObjects.Item(1)
$t.Resize($ws.Range("B5:G96"))

# Update the "last updated" caption to reflect the new month. (This text
# lived in B96 before the row insert; after inserting at row 6 it shifted
# down to B97.)
$ws.Range("B97").Value = "Actualización: Julio 2025."
